$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (old Neutrophils sending-cluster rows removed entirely)
$ws.Rows("4:5").Delete()

# Row 2: keep A2="ECs", B2="Pomc", C2="Oprm1"; change D2 target cluster to "Neutrophils"
$ws.Range("D2").Value = "Neutrophils"

# Row 2 numeric metric updates
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.8055585
$ws.Range("H2").Value = 1.611117
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.896219333333333
$ws.Range("N2").Value = 5.688658
$ws.Range("O2").Value = 0.9918549325388253
$ws.Range("P2").Value = 0.9918549325388253
$ws.Range("Q2").Value = 1.527515601831
$ws.Range("R2").Value = 9.165093610986
$ws.Range("S2").Value = 0.9918549325388253
$ws.Range("T2").Value = 0.9918549325388253

# Row 3 keeps A3="ECs", B3="Pomc", C3="Oprm1", D3="Resolving-Mac" (unchanged)
# Row 3 numeric metric updates
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.8055585
$ws.Range("H3").Value = 1.611117
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01557166666666667
$ws.Range("N3").Value = 0.046715
$ws.Range("O3").Value = 0.008145067461174712
$ws.Range("P3").Value = 0.008145067461174714
$ws.Range("Q3").Value = 0.0125438884425
$ws.Range("R3").Value = 0.075263330655
$ws.Range("S3").Value = 0.008145067461174712
$ws.Range("T3").Value = 0.008145067461174714
